$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Price) and E (Volume 1h) to be treated as text so that
# numeric-looking strings (e.g. "310.50", "-3.39%") are preserved exactly
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "310.50"
$ws.Cells.Item(2, 5).Value = "-3.39%"
$ws.Cells.Item(3, 4).Value = "54.30"
$ws.Cells.Item(3, 5).Value = "10.48%"
$ws.Cells.Item(4, 4).Value = "5.095"
$ws.Cells.Item(4, 5).Value = "-4.53%"
$ws.Cells.Item(5, 4).Value = "0.07930"
$ws.Cells.Item(5, 5).Value = "-1.77%"
$ws.Cells.Item(6, 4).Value = "4.551"
$ws.Cells.Item(6, 5).Value = "-0.97%"
$ws.Cells.Item(7, 5).Value = "4.24%"
$ws.Cells.Item(8, 4).Value = "1.675"
$ws.Cells.Item(8, 5).Value = "1.86%"
$ws.Cells.Item(9, 4).Value = "0.1245"
$ws.Cells.Item(9, 5).Value = "-2.85%"
$ws.Cells.Item(10, 4).Value = "0.2011"
$ws.Cells.Item(10, 5).Value = "2.27%"
$ws.Cells.Item(11, 4).Value = "0.04733"
$ws.Cells.Item(11, 5).Value = "1.04%"
$ws.Cells.Item(12, 4).Value = "0.09404"
$ws.Cells.Item(12, 5).Value = "-3.81%"
$ws.Cells.Item(13, 4).Value = "0.1042"
$ws.Cells.Item(13, 5).Value = "-0.42%"
$ws.Cells.Item(14, 5).Value = "-3.49%"
$ws.Cells.Item(15, 4).Value = "0.005809"
$ws.Cells.Item(15, 5).Value = "-0.38%"
$ws.Cells.Item(16, 5).Value = "2,015.12%"
$ws.Cells.Item(17, 4).Value = "3.336"
$ws.Cells.Item(17, 5).Value = "-0.17%"
$ws.Cells.Item(18, 4).Value = "2.436"
$ws.Cells.Item(18, 5).Value = "-0.43%"
$ws.Cells.Item(19, 4).Value = "0.3429"
$ws.Cells.Item(19, 5).Value = "-2.21%"
$ws.Cells.Item(20, 4).Value = "8.332"
$ws.Cells.Item(20, 5).Value = "4.25%"
$ws.Cells.Item(21, 4).Value = "0.1360"
$ws.Cells.Item(21, 5).Value = "-0.45%"
$ws.Cells.Item(22, 4).Value = "0.2903"
$ws.Cells.Item(22, 5).Value = "-6.14%"
$ws.Cells.Item(23, 4).Value = "0.04160"
$ws.Cells.Item(23, 5).Value = "-0.86%"
$ws.Cells.Item(24, 4).Value = "0.001251"
$ws.Cells.Item(24, 5).Value = "-4.60%"
$ws.Cells.Item(25, 4).Value = "0.003974"
$ws.Cells.Item(25, 5).Value = "-6.97%"
$ws.Cells.Item(26, 4).Value = "0.0001341"
$ws.Cells.Item(26, 5).Value = "-0.50%"
$ws.Cells.Item(38, 4).Value = "0.02640"
$ws.Cells.Item(38, 5).Value = "-3.08%"
$ws.Cells.Item(39, 4).Value = "0.05948"
$ws.Cells.Item(39, 5).Value = "-2.97%"
$ws.Cells.Item(40, 4).Value = "0.01083"
$ws.Cells.Item(40, 5).Value = "-0.11%"
$ws.Cells.Item(41, 5).Value = "19.19%"
$ws.Cells.Item(42, 4).Value = "0.007934"
$ws.Cells.Item(42, 5).Value = "-1.01%"
$ws.Cells.Item(43, 4).Value = "0.008137"
$ws.Cells.Item(43, 5).Value = "3.04%"
$ws.Cells.Item(44, 4).Value = "0.008338"
$ws.Cells.Item(44, 5).Value = "-3.62%"
$ws.Cells.Item(45, 4).Value = "0.3390"
$ws.Cells.Item(45, 5).Value = "-3.16%"
$ws.Cells.Item(46, 4).Value = "0.00007256"
$ws.Cells.Item(46, 5).Value = "5.84%"
$ws.Cells.Item(47, 4).Value = "0.00000000746"
$ws.Cells.Item(47, 5).Value = "-0.35%"
$ws.Cells.Item(48, 4).Value = "0.05687"
$ws.Cells.Item(48, 5).Value = "2.66%"
$ws.Cells.Item(49, 4).Value = "0.002605"
$ws.Cells.Item(49, 5).Value = "-34.87%"
$ws.Cells.Item(50, 4).Value = "0.00002089"
$ws.Cells.Item(50, 5).Value = "-0.35%"
$ws.Cells.Item(51, 4).Value = "0.0001989"
$ws.Cells.Item(51, 5).Value = "-0.35%"
